$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHA_RETAU5200")

# --- 1. Fill in O89 (second-order strain-rate approximation at the first layer) ---
$ws.Range("O89").Value = 5163.6548122326203

# --- 2. Add four new rows (92:95) for the WM+SMAG+SMALL+ONESIDE case ---
# Prime formatting by copying the analogous block (rows 86:89) down to 92:95
$ws.Range("A86:X89").Copy()
$ws.Range("A92:X95").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$rows = @(92, 93, 94, 95)
$dvals = @("128×48×32", "192×72×48", "256×96×64", "384×144×96")
$evals = @(0.1, 0.067, 0.05, 0.033)
$ovals = @(5161.2796463928098, 5137.6418223186201, 5175.4405682218103, 5160.3492684698103)

for ($i = 0; $i -lt 4; $i++) {
    $r = $rows[$i]
    $ws.Range("A$r").Value = "WM+SMAG+SMALL+ONESIDE"
    $ws.Range("B$r").Value = 250000
    $ws.Range("C$r").Value = "12.8×4.8×2.0"
    $ws.Range("D$r").Value = $dvals[$i]
    $ws.Range("E$r").Value = $evals[$i]
    $ws.Range("F$r").Value = $evals[$i]
    $ws.Range("G$r").Formula = "=F$r"
    $ws.Range("H$r").Formula = "=0.25*G$r"
    $ws.Range("I$r").Value = 1
    $ws.Range("J$r").Value = "stretching"
    $ws.Range("K$r").Value = "CFR"
    $ws.Range("L$r").Value = 0.1
    $ws.Range("M$r").Formula = "=L$r/H$r"
    $ws.Range("N$r").Value = 5185.8969999999999
    $ws.Range("O$r").Value = $ovals[$i]
    $ws.Range("P$r").Formula = "=8*(N$r/B$r)^2"
}

# Column Q: Q92 and Q93 are entered individually; Q94:Q95 share one formula group.
$ws.Range("Q92").Formula = "=8*(O92/B92)^2"
$ws.Range("R92").Formula = "=(Q92-P92)/P92"
$ws.Range("S92").Formula = "=500*2*O92/B92"
$ws.Range("T92").Formula = "=B92/4*P92"

# Columns U,V,W,X share one formula group each across the full 92:95 block.
$ws.Range("U92:U95").Formula = "=E92*N92"
$ws.Range("V92:V95").Formula = "=F92*O92"
$ws.Range("W92:W95").Formula = "=G92*N92"
$ws.Range("X92:X95").Formula = "=H92*N92"

$ws.Range("Q93").Formula = "=8*(O93/B93)^2"

# Columns R,S,T share one formula group each across the 93:95 block.
$ws.Range("R93:R95").Formula = "=(Q93-P93)/P93"
$ws.Range("S93:S95").Formula = "=500*2*O93/B93"
$ws.Range("T93:T95").Formula = "=B93/4*P93"

# Column Q: remaining 94:95 share one formula group.
$ws.Range("Q94:Q95").Formula = "=8*(O94/B94)^2"

# --- 3. Update the sheet view (scroll position + active selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("R92").Select()
